$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 7.239
$ws.Range("A8").Value = -21.672
$ws.Range("A10").Value = -21.794
$ws.Range("B11").Value = 6.776999999999999
$ws.Range("A12").Value = -21.53
$ws.Range("B12").Value = 5.397
$ws.Range("B15").Value = 5.776
$ws.Range("B17").Value = 4.994000000000001
$ws.Range("A18").Value = -21.885
$ws.Range("A25").Value = -21.75
$ws.Range("B26").Value = 6.319000000000001
$ws.Range("B27").Value = 5.994999999999999
$ws.Range("B28").Value = 5.815
$ws.Range("B32").Value = 6.859
$ws.Range("A37").Value = -21.081
$ws.Range("B37").Value = 7.888000000000001
$ws.Range("B41").Value = 7.773000000000001
$ws.Range("B47").Value = 5.69
$ws.Range("B51").Value = 6.664
$ws.Range("A55").Value = -21.843
$ws.Range("B65").Value = 6.282000000000001
$ws.Range("A68").Value = -21.583
$ws.Range("B73").Value = 6.528
$ws.Range("A77").Value = -21.098
$ws.Range("A78").Value = -20.739
$ws.Range("A79").Value = -21.233
$ws.Range("A80").Value = -20.887
$ws.Range("A81").Value = -21.63
$ws.Range("A82").Value = -21.886
$ws.Range("A84").Value = -20.552
$ws.Range("B84").Value = 7.964999999999999
$ws.Range("B85").Value = 6.083
$ws.Range("B89").Value = 5.205000000000001
$ws.Range("B93").Value = 6.077000000000001
$ws.Range("B95").Value = 5.898999999999999
$ws.Range("B98").Value = 7.1
$ws.Range("B99").Value = 5.840000000000001
$ws.Range("A101").Value = -22.156
$ws.Range("B101").Value = 5.328
$ws.Range("A102").Value = -21.291
$ws.Range("B102").Value = 6.483
